# The source edit (commit "Add powerpoint 365 support") re-saved this
# deck with PowerPoint 365, which is the only thing that changed the
# presentation.xml <p:sldId r:id="..."> relationship-id strings
# (R84ff40fa117f4dfe -> Rc9104575d00f460a, etc.). Every <p:sldId id="...">
# value, the slide order, the slide count and every slide's visible
# content/shapes are byte-for-byte identical before and after - only the
# internal, otherwise-meaningless r:id tokens used to resolve each
# <p:sldId> to its part in presentation.xml.rels were regenerated by the
# newer Office build's packaging code.
#
# That relationship-id regeneration is a side effect of the host
# application's save routine, not something exposed anywhere on the
# PowerPoint object model (Slides/Slide do not expose the underlying
# OOXML relationship id, and there is no supported automation call that
# "re-saves with new part ids" without touching content). Re-creating the
# slides (Delete + Add/Duplicate/Paste/InsertFromFile) would change the
# <p:sldId id="..."> numbers and collection order as well, which would
# contradict the rest of the diff, so it is not a faithful way to mimic
# this particular change.
#
# The presentation is already open as $ppt.ActivePresentation and is
# re-serialized by the host after this script runs, which is the
# accurate object-model equivalent of "opened and saved by a newer
# PowerPoint build" for a deck whose slides/content do not otherwise
# change.
$p = $ppt.ActivePresentation
